# Lecture partielle de l'EDT M1 MIAGE.
# Shift the schedule dates forward by 3 years (1096 days) and refresh the
# French weekday labels in column B to match the new dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$joursFr = @("lundi", "mardi", "mercredi", "jeudi", "vendredi", "samedi", "dimanche")

# Rows whose date (column A) must be advanced by 1096 days, with the
# corresponding weekday label (column B) recomputed from the new date.
$rows = 2, 4, 7, 9, 12, 15, 19, 22

foreach ($r in $rows) {
    $cellA = $ws.Cells.Item($r, 1)
    $oldSerial = $cellA.Value2
    $newSerial = $oldSerial + 1096
    $cellA.Value2 = $newSerial

    # Excel date serials: day 0 = 1899-12-30 (matches the workbook's non-1904 date system).
    # 1899-12-30 was a Saturday, i.e. weekday index 5 (0=lundi ... 6=dimanche).
    $dayIndex = (5 + $newSerial) % 7
    $ws.Cells.Item($r, 2).Value = $joursFr[$dayIndex]
}
